# "updating for summer 2019": change the cover-page term/year and move the
# editing cursor's "_GoBack" bookmark to where Word last left it.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Cover page: "Spring 2019" (spread across 3 runs: "Spring" / " 201" /
#    "9") -> "Fall 2018" (single run). A plain Find/Replace over that
#    paragraph mirrors what a user would do and naturally collapses the
#    old runs into one, inheriting the first run's formatting.
# ---------------------------------------------------------------------
$titleRange = $d.Paragraphs.Item(2).Range
$titleRange.Find.Execute("Spring 2019", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "Fall 2018", 2)

# ---------------------------------------------------------------------
# 2. Word keeps a hidden "_GoBack" bookmark at the location of the most
#    recent edit. It used to sit right after the title text; now that the
#    title edit happened, it belongs at the spot the author worked on
#    last instead (the blank Book Title paragraph right before the
#    "Introduction" section, just after the Table of Contents).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$lastEditPara = $d.Paragraphs.Item(24)
$lastEditPara.Range.Select()
$d.Bookmarks.Add("_GoBack", $word.Selection.Range)
